$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting existing rows 23:109 down to 24:110.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new data record.
$ws.Range("A23").Value2 = 4
$ws.Range("B23").Value2 = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C23").Value2 = 'Los Lagos'
$ws.Range("D23").Value2 = 44600
$ws.Range("E23").Value2 = 10
$ws.Range("F23").Value2 = 100112022
$ws.Range("G23").Value2 = 'Arveja Verde'
$ws.Range("H23").Value2 = 'Sin especificar'
$ws.Range("I23").Value2 = 'Primera'
$ws.Range("J23").Value2 = 80
$ws.Range("K23").Value2 = 23000
$ws.Range("L23").Value2 = 23000
$ws.Range("M23").Value2 = 23000
$ws.Range("N23").Value2 = '$/saco 25 kilos'
$ws.Range("O23").Value2 = 'Región de La Araucanía'
$ws.Range("P23").Value2 = 920
$ws.Range("Q23").Value2 = 25
$ws.Range("R23").Value2 = 'Hortaliza'
